$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First paragraph: pad the existing sentence with two trailing spaces and
#    append a new red "(This is a change - Version for branch alternate)"
#    annotation, split across three separate runs (all colored C00000).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$null = $p1.Range.Find.Execute(
    "This is a Microsoft word document.", $true, $false, $false, $false,
    $false, $true, 1, $false, "This is a Microsoft word document.  ", 2)

$redColor = 192  # wdColor value for RRGGBB = C00000 (stored as 0x00BBGGRR)

$parts = @(
    "(This is a change " + [char]0x2013 + " Ve",
    "rsion for branch alternate",
    ")"
)

foreach ($part in $parts) {
    $p1 = $d.Paragraphs(1)
    $r = $p1.Range
    $startPos = $r.End
    $r.InsertAfter($part)
    $endPos = $r.End
    # $r.End (both before/after) sits one past the paragraph mark, so shift
    # back by one to address only the text that was just inserted.
    $newRange = $d.Range($startPos - 1, $endPos - 1)
    $newRange.Font.Color = $redColor
}

# ---------------------------------------------------------------------------
# 2) Append a new, empty paragraph after the final paragraph, shaded with
#    fill color F9F9F9 and no other paragraph/run formatting.
# ---------------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$lastPara.Range.InsertParagraphAfter()

$newIndex = $d.Paragraphs.Count
$newPara = $d.Paragraphs($newIndex)
$newPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>')

Write-Output "done"
